$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "ScoreF" (2nd tab) - add J/K/L (Minute3/Second3/Rep3) columns
# ------------------------------------------------------------------
$wsF = $wb.Worksheets.Item(2)

$scoreF = @{
    2  = @(13, 0, 232)
    3  = @(13, 0, 238)
    4  = @(11, 26, 258)
    5  = @(10, 13, 258)
    6  = @(11, 34, 258)
    7  = @(12, 47, 258)
    8  = @(10, 31, 258)
    9  = @(13, 0, 237)
    10 = @(13, 0, 157)
    11 = @(12, 47, 258)
    12 = @(12, 56, 258)
    13 = @(13, 0, 237)
    14 = @(13, 0, 122)
    15 = @(13, 0, 216)
    16 = @(11, 29, 258)
    17 = @(13, 0, 1)
    18 = @(13, 0, 91)
    19 = @(8, 24, 258)
    20 = @(10, 21, 258)
    21 = @(10, 33, 258)
    22 = @(13, 0, 195)
    23 = @(12, 36, 258)
    24 = @(12, 43, 258)
    25 = @(13, 0, 226)
    26 = @(13, 0, 254)
}

foreach ($row in $scoreF.Keys) {
    $vals = $scoreF[$row]
    $wsF.Cells.Item($row, 10).Value2 = $vals[0]
    $wsF.Cells.Item($row, 11).Value2 = $vals[1]
    $wsF.Cells.Item($row, 12).Value2 = $vals[2]
}

# ------------------------------------------------------------------
# Sheet "SFM" (3rd tab) - add a ranking table (rows 2-9) with widened
# columns to match the new data
# ------------------------------------------------------------------
$wsM = $wb.Worksheets.Item(3)

$matrix = @(
    @("Ole og Mikus", "Fit and Lazy", 20),
    @("Henrik og Henrik", "OnlyFans", 18),
    @("Jakob og Finn", "The NHH Nerds", 16),
    @("Aril og Mats", "Maverick & Goose", 14),
    @("Mathias og Rakan", "Rakus Mathius", 12),
    @("Jonas og Arild", "Team Kongobajer", 10),
    @("Thor Andre og Ole Andre", "Skaol Skadle", 8),
    @("Morten og Gabor", "Daddszkys", 6)
)

$r = 2
foreach ($entry in $matrix) {
    $wsM.Cells.Item($r, 1).Value2 = $entry[0]
    $wsM.Cells.Item($r, 2).Value2 = $entry[1]
    $wsM.Cells.Item($r, 3).Value2 = $entry[2]
    $r = $r + 1
}

# (ColumnWidth is internally quantized to 1/6-character steps by the
# runtime, so these inputs are chosen to land as close as possible to
# the target "best fit" widths of 21.109375 / 7.77734375 / 7.88671875 /
# 5.109375 characters.)
$wsM.Columns.Item(1).ColumnWidth = 20.333333333333336
$wsM.Columns.Item(4).ColumnWidth = 7
$wsM.Columns.Item(5).ColumnWidth = 7
$wsM.Columns.Item(6).ColumnWidth = 4.333333333333334
$wsM.Columns.Item(7).ColumnWidth = 7
$wsM.Columns.Item(8).ColumnWidth = 7
$wsM.Columns.Item(9).ColumnWidth = 4.333333333333334

# ------------------------------------------------------------------
# Sheet selections / active tab:
#  - ScoreM (1st tab) loses tabSelected, selection becomes A2:B23
#  - ScoreF (2nd tab) becomes the selected/active tab, selection K24
#  - SFM (3rd tab) selection becomes C10
# ------------------------------------------------------------------
$wsM.Activate()
$wsM.Range("C10").Select()

$wsScoreM = $wb.Worksheets.Item(1)
$wsScoreM.Activate()
$wsScoreM.Range("A2:B23").Select()

$wsF.Activate()
$wsF.Range("K24").Select()
